$d = $word.ActiveDocument

$d.Content.Find.Execute("718×9=6462", $true, $false, $false, $false, $false, $true, 1, $false, "148×6=888", 2) | Out-Null
$d.Content.Find.Execute("633×7=4431", $true, $false, $false, $false, $false, $true, 1, $false, "915×5=4575", 2) | Out-Null
$d.Content.Find.Execute("170×4=680", $true, $false, $false, $false, $false, $true, 1, $false, "874×7=6118", 2) | Out-Null
$d.Content.Find.Execute("853×9=7677", $true, $false, $false, $false, $false, $true, 1, $false, "154×2=308", 2) | Out-Null
$d.Content.Find.Execute("474×8=3792", $true, $false, $false, $false, $false, $true, 1, $false, "563×2=1126", 2) | Out-Null
$d.Content.Find.Execute("532×4=2128", $true, $false, $false, $false, $false, $true, 1, $false, "236×4=944", 2) | Out-Null
$d.Content.Find.Execute("285×2=570", $true, $false, $false, $false, $false, $true, 1, $false, "796×3=2388", 2) | Out-Null
$d.Content.Find.Execute("432×9=3888", $true, $false, $false, $false, $false, $true, 1, $false, "504×2=1008", 2) | Out-Null
$d.Content.Find.Execute("947×8=7576", $true, $false, $false, $false, $false, $true, 1, $false, "892×7=6244", 2) | Out-Null
$d.Content.Find.Execute("805×9=7245", $true, $false, $false, $false, $false, $true, 1, $false, "268×6=1608", 2) | Out-Null
$d.Content.Find.Execute("467×6=2802", $true, $false, $false, $false, $false, $true, 1, $false, "228×2=456", 2) | Out-Null
$d.Content.Find.Execute("916×4=3664", $true, $false, $false, $false, $false, $true, 1, $false, "393×3=1179", 2) | Out-Null
$d.Content.Find.Execute("284×6=1704", $true, $false, $false, $false, $false, $true, 1, $false, "543×7=3801", 2) | Out-Null
$d.Content.Find.Execute("161×9=1449", $true, $false, $false, $false, $false, $true, 1, $false, "509×6=3054", 2) | Out-Null
$d.Content.Find.Execute("483×3=1449", $true, $false, $false, $false, $false, $true, 1, $false, "204×6=1224", 2) | Out-Null
$d.Content.Find.Execute("307×5=1535", $true, $false, $false, $false, $false, $true, 1, $false, "466×6=2796", 2) | Out-Null
$d.Content.Find.Execute("193×8=1544", $true, $false, $false, $false, $false, $true, 1, $false, "840×4=3360", 2) | Out-Null
$d.Content.Find.Execute("732×3=2196", $true, $false, $false, $false, $false, $true, 1, $false, "598×4=2392", 2) | Out-Null
$d.Content.Find.Execute("402×6=2412", $true, $false, $false, $false, $false, $true, 1, $false, "807×7=5649", 2) | Out-Null
$d.Content.Find.Execute("849×6=5094", $true, $false, $false, $false, $false, $true, 1, $false, "131×7=917", 2) | Out-Null
$d.Content.Find.Execute("523×9=4707", $true, $false, $false, $false, $false, $true, 1, $false, "924×8=7392", 2) | Out-Null
$d.Content.Find.Execute("503×3=1509", $true, $false, $false, $false, $false, $true, 1, $false, "632×3=1896", 2) | Out-Null
$d.Content.Find.Execute("240×9=2160", $true, $false, $false, $false, $false, $true, 1, $false, "835×2=1670", 2) | Out-Null
$d.Content.Find.Execute("603×5=3015", $true, $false, $false, $false, $false, $true, 1, $false, "659×5=3295", 2) | Out-Null
$d.Content.Find.Execute("276×2=552", $true, $false, $false, $false, $false, $true, 1, $false, "966×6=5796", 2) | Out-Null
